# paises.xlsx data refresh (COVID dashboard) + provincias Spain.
# - Row 1: "Datos actualizados..." timestamp bumped from 13:10 to 14:27.
# - Four pairs/groups of countries swapped row-adjacent positions because the
#   underlying table is sorted by total cases (col B) descending and the new
#   day's numbers changed the ranking: Costa Rica/Nepal (rows 55-56),
#   Guinea-Bisau/Benin/Islandia (rows 148-150), Santa Lucia/Timor Oriental
#   (rows 204-205) and Montserrat/Islas Malvinas (rows 214-215).
# - Several rows' B-H metrics (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes hoy, Muertes) were refreshed with the
#   latest counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reshuffles (sharedStrings churn caused by re-sort) ---
$ws.Cells.Item(55, 1).Value = "Nepal"
$ws.Cells.Item(56, 1).Value = "Costa Rica"
$ws.Cells.Item(148, 1).Value = "Islandia"
$ws.Cells.Item(149, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(150, 1).Value = "Benin"
$ws.Cells.Item(204, 1).Value = "Timor Oriental"
$ws.Cells.Item(205, 1).Value = "Santa Lucia"
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 1).Value = "Montserrat"

# --- Timestamp update ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 19 de Septiembre de 2020 a las 14:27"

# --- Refreshed numeric data (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes) ---
$ws.Cells.Item(4, 2).Value = 6928304
$ws.Cells.Item(4, 3).Value = 2363
$ws.Cells.Item(4, 4).Value = 4192774
$ws.Cells.Item(4, 5).Value = 2532342
$ws.Cells.Item(4, 7).Value = 17
$ws.Cells.Item(4, 8).Value = 203188
$ws.Cells.Item(13, 4).Value = 478077
$ws.Cells.Item(13, 5).Value = 122876
$ws.Cells.Item(13, 7).Value = 49
$ws.Cells.Item(13, 8).Value = 12705
$ws.Cells.Item(32, 2).Value = 123146
$ws.Cells.Item(32, 3).Value = 229
$ws.Cells.Item(32, 4).Value = 120089
$ws.Cells.Item(32, 5).Value = 2848
$ws.Cells.Item(46, 2).Value = 84242
$ws.Cells.Item(46, 3).Value = 809
$ws.Cells.Item(46, 4).Value = 73512
$ws.Cells.Item(46, 5).Value = 10326
$ws.Cells.Item(46, 7).Value = 1
$ws.Cells.Item(46, 8).Value = 404
$ws.Cells.Item(55, 2).Value = 62797
$ws.Cells.Item(55, 3).Value = 1204
$ws.Cells.Item(55, 4).Value = 45267
$ws.Cells.Item(55, 5).Value = 17129
$ws.Cells.Item(55, 7).Value = 11
$ws.Cells.Item(55, 8).Value = 401
$ws.Cells.Item(56, 2).Value = 62374
$ws.Cells.Item(56, 4).Value = 23160
$ws.Cells.Item(56, 5).Value = 38528
$ws.Cells.Item(56, 8).Value = 686
$ws.Cells.Item(87, 2).Value = 16020
$ws.Cells.Item(87, 3).Value = 49
$ws.Cells.Item(87, 4).Value = 14630
$ws.Cells.Item(87, 5).Value = 1171
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 219
$ws.Cells.Item(90, 2).Value = 14688
$ws.Cells.Item(90, 3).Value = 43
$ws.Cells.Item(90, 4).Value = 11153
$ws.Cells.Item(90, 5).Value = 3233
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = 302
$ws.Cells.Item(103, 2).Value = 8922
$ws.Cells.Item(103, 3).Value = 64
$ws.Cells.Item(103, 5).Value = 883
$ws.Cells.Item(125, 5).Value = 3410
$ws.Cells.Item(125, 7).Value = 11
$ws.Cells.Item(125, 8).Value = 81
$ws.Cells.Item(148, 2).Value = 2307
$ws.Cells.Item(148, 3).Value = 77
$ws.Cells.Item(148, 4).Value = 2116
$ws.Cells.Item(148, 5).Value = 181
$ws.Cells.Item(148, 8).Value = 10
$ws.Cells.Item(149, 2).Value = 2303
$ws.Cells.Item(149, 4).Value = 1127
$ws.Cells.Item(149, 5).Value = 1137
$ws.Cells.Item(149, 8).Value = 39
$ws.Cells.Item(150, 2).Value = 2280
$ws.Cells.Item(150, 4).Value = 1950
$ws.Cells.Item(150, 5).Value = 290
$ws.Cells.Item(150, 8).Value = 40
$ws.Cells.Item(168, 4).Value = 942
$ws.Cells.Item(168, 5).Value = 91
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1
